# Fixing issues with experiments validation
# Target sheet: "chip-seq input dna" (5th tab)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chip-seq input dna")

# Insert 3 new columns before column B (shifts existing B..AA right to E..AD)
$ws.Range("B:D").Insert()

# New header row values for the inserted columns
$ws.Range("B1").Value = "Experiment Alias"
$ws.Range("C1").Value = "Project"
$ws.Range("D1").Value = "Secondary Project"

# Match the column widths Excel's "best fit" produced for the new headers
$ws.Columns.Item(2).ColumnWidth = 14.25
$ws.Columns.Item(3).ColumnWidth = 5.92
$ws.Columns.Item(4).ColumnWidth = 14.92

# Move the active selection to A2, as in the edited workbook
$ws.Range("A2").Select()

Write-Host "Applied chip-seq input dna header edits"
